$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated FBA results after removing the minYield constraint from the solver;
# rateFBA (E), yieldFBA (F), vCarbonSourcesFBA (G) and the yRBA/yFBA ratio (K)
# were recomputed per product row. RBA columns (H/I/J) are untouched.
$ws.Range("E2").Value = 17.81382210473816
$ws.Range("F2").Value = 0.5695516925744543
$ws.Range("K2").Value = 0.5800987514688722

$ws.Range("E3").Value = 6.941370044788568
$ws.Range("F3").Value = 0.3563056204478834
$ws.Range("G3").Value = -13.21
$ws.Range("K3").Value = 0.5425276625108567

$ws.Range("E4").Value = 20.06614443981965
$ws.Range("F4").Value = 0.7595227088575357
$ws.Range("G4").Value = -13.21
$ws.Range("K4").Value = 0.2387059039238457

$ws.Range("E5").Value = 22.34089249430825
$ws.Range("F5").Value = 0.8456240926831734
$ws.Range("G5").Value = -13.21
$ws.Range("K5").Value = 0.3703426488291691

$ws.Range("E6").Value = 20.06614443981965
$ws.Range("F6").Value = 0.7595227088575357
$ws.Range("G6").Value = -13.21
$ws.Range("K6").Value = 0.2387059039238457

$ws.Range("E7").Value = 22.34089249430825
$ws.Range("F7").Value = 0.8456240926831734
$ws.Range("G7").Value = -13.21
$ws.Range("K7").Value = 0.3703426488291691

$ws.Range("E8").Value = 9.055951924829568
$ws.Range("F8").Value = 0.521737065160645
$ws.Range("G8").Value = -13.21
$ws.Range("K8").Value = 0.3535611162697625

$ws.Range("E9").Value = 2.844697670245594
$ws.Range("F9").Value = 0.2800995824094558
$ws.Range("K9").Value = 0.4459930778203139

$ws.Range("E10").Value = 12.8944378446668
$ws.Range("F10").Value = 0.488283498804935
$ws.Range("G10").Value = -13.21
$ws.Range("K10").Value = 0.9974082223302524

$ws.Range("E11").Value = 11.19787294998664
$ws.Range("F11").Value = 0.3487541611690804
$ws.Range("G11").Value = -13.21
$ws.Range("K11").Value = 0.7847593657794829

$ws.Range("E12").Value = 12.48937108880209
$ws.Range("F12").Value = 0.7667243171293706
$ws.Range("K12").Value = 0.3990386909380163

$ws.Range("E13").Value = 7.096204221127857
$ws.Range("F13").Value = 0.4864762418608935
$ws.Range("G13").Value = -13.21
$ws.Range("K13").Value = 0.4351111423546155

$ws.Range("E14").Value = 2.761960154328127
$ws.Range("F14").Value = 0.2975663935774678
$ws.Range("G14").Value = -13.21
$ws.Range("K14").Value = 0.2594030863872557

$ws.Range("E15").Value = 2.677177417204605
$ws.Range("F15").Value = 0.2727280956668592
$ws.Range("G15").Value = -13.21
$ws.Range("K15").Value = 0.2720593940343134

$ws.Range("E16").Value = 12.5238538026001
$ws.Range("F16").Value = 0.3900514094987304
$ws.Range("G16").Value = -13.21
$ws.Range("K16").Value = 0.9872554154549937

$ws.Range("E17").Value = 24.24659146894868
$ws.Range("F17").Value = 0.9074664948630266
$ws.Range("G17").Value = -13.21
$ws.Range("K17").Value = 0.8724089016032093

$ws.Range("E18").Value = 10.42241532098108
$ws.Range("F18").Value = 0.6135132722286827
$ws.Range("G18").Value = -13.21
$ws.Range("K18").Value = 0.4612523037661422

$ws.Range("E19").Value = 3.969093061921526
$ws.Range("F19").Value = 0.4540541212621765
$ws.Range("K19").Value = 0.2713272712592787

$ws.Range("E20").Value = 0.09996363266789167
$ws.Range("F20").Value = 0.2355252527240741
$ws.Range("K20").Value = 1.008269443564308

$ws.Range("E21").Value = 0.114167371871854
$ws.Range("F21").Value = 0.2428831784355429
$ws.Range("K21").Value = 1.007992630508971

$ws.Range("E22").Value = 12.48937108880229
$ws.Range("F22").Value = 0.5410628137990376
$ws.Range("G22").Value = -13.21
$ws.Range("K22").Value = 0.3818254643828795

$ws.Range("E23").Value = 2.511143887422177
$ws.Range("F23").Value = 0.3475710892081204
$ws.Range("K23").Value = 0.1959094385288666

$ws.Range("E24").Value = 3.96909303424228
$ws.Range("F24").Value = 0.3806549565258843
$ws.Range("G24").Value = -13.21
$ws.Range("K24").Value = 0.2713282279411472

$ws.Range("E25").Value = 2.957467015229745
$ws.Range("F25").Value = 0.3060784639564738
$ws.Range("G25").Value = -13.21
$ws.Range("K25").Value = 0.4476413612563981

$ws.Range("E26").Value = 9.725652288380664
$ws.Range("F26").Value = 0.7116896813834993
$ws.Range("K26").Value = 0.4692139917052498

$ws.Range("E27").Value = 0.1013955785655886
$ws.Range("F27").Value = 0.2406785828175028
$ws.Range("K27").Value = 1.008402006346098

$ws.Range("E28").Value = 6.766253948281869
$ws.Range("F28").Value = 0.2961120081300709
$ws.Range("G28").Value = -13.21
$ws.Range("K28").Value = 0.4418733056202847

$ws.Range("E29").Value = 0.1308423118282197
$ws.Range("F29").Value = 0.2719810456277346
$ws.Range("K29").Value = 0.4515292272393641

$ws.Range("E30").Value = 0.007105611207608455
$ws.Range("F30").Value = 0.06718413215330893

